# Bop.Calc.xlsx edit: rename the "index" column to "i" and make it
# zero-based instead of one-based (refactor rename index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the table/header column from "index" to "i" (renaming the header
# cell automatically renames the corresponding ListColumn of the
# "testdata" table that spans A1:I503).
$ws.Range("A1").Value = "i"

# The column held a 1-based row counter (1..502). Make it 0-based (0..501)
# by decrementing every data value by one.
$dataRange = $ws.Range("A2:A503")
$values = $dataRange.Value()
for ($r = 1; $r -le 502; $r++) {
    $values[$r, 1] = $values[$r, 1] - 1
}
$dataRange.Value = $values

# Narrow column A now that it only needs to fit a shorter number/heading
# (stored column width of 6 -> 4 characters).
$ws.Columns.Item(1).ColumnWidth = 3.1666666666666665
